$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 160, pushing the existing rows 160-163
# (two weekly records) down to 162-165, then fill the two new rows with
# the latest week's data (fecha 2021-11-09 / serial 44509).
$ws.Rows.Item(160).Resize(2).Insert()

# New row 160: Primera quality record for the newest week
$ws.Cells.Item(160, 1).Value = 11
$ws.Cells.Item(160, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(160, 3).Value = "Bíobío"
$ws.Cells.Item(160, 4).Value = 44509
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(159, 4).NumberFormat
$ws.Cells.Item(160, 5).Value = 8
$ws.Cells.Item(160, 6).Value = 100112017
$ws.Cells.Item(160, 7).Value = "Apio"
$ws.Cells.Item(160, 8).Value = "Americana (o)"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 100
$ws.Cells.Item(160, 11).Value = 7500
$ws.Cells.Item(160, 12).Value = 8000
$ws.Cells.Item(160, 13).Value = 7750
$ws.Cells.Item(160, 14).Value = "`$/docena de matas"
$ws.Cells.Item(160, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(160, 16).Value = 1292
$ws.Cells.Item(160, 17).Value = 6
$ws.Cells.Item(160, 18).Value = "Hortaliza"

# New row 161: Segunda quality record for the newest week
$ws.Cells.Item(161, 1).Value = 11
$ws.Cells.Item(161, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(161, 3).Value = "Bíobío"
$ws.Cells.Item(161, 4).Value = 44509
$ws.Cells.Item(161, 4).NumberFormat = $ws.Cells.Item(159, 4).NumberFormat
$ws.Cells.Item(161, 5).Value = 8
$ws.Cells.Item(161, 6).Value = 100112017
$ws.Cells.Item(161, 7).Value = "Apio"
$ws.Cells.Item(161, 8).Value = "Americana (o)"
$ws.Cells.Item(161, 9).Value = "Segunda"
$ws.Cells.Item(161, 10).Value = 50
$ws.Cells.Item(161, 11).Value = 6500
$ws.Cells.Item(161, 12).Value = 6500
$ws.Cells.Item(161, 13).Value = 6500
$ws.Cells.Item(161, 14).Value = "`$/docena de matas"
$ws.Cells.Item(161, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(161, 16).Value = 1083
$ws.Cells.Item(161, 17).Value = 6
$ws.Cells.Item(161, 18).Value = "Hortaliza"
